$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 13, shifting the old rows 13-23 down to 14-24.
# (This fixes a pre-existing off-by-one misalignment between the A-column labels
#  and their B/C content, and makes room for the newly authored cell text below.)
$ws.Rows.Item(13).Insert()

# Fill in the B/C content that the row insert left blank / mis-stocked with the
# new authored text for this row.
$ws.Range("B10").Value = 'Visão integrada sobre petróleo e gás natural, desde a origem até o processamento primário. Descrições, características e aplicações dos derivados do petróleo. Processo e esquemas de refino e processamento do gás natural.'
$ws.Range("C10").Value = 'Visão integrada sobre petróleo e gás natural, desde a origem até o processamento primário. Descrições, características e aplicações dos derivados do petróleo. Processo e esquemas de refino e processamento do gás natural.'

$ws.Range("B13").Value = '1285870 - Marcos Villela Barcza'
$ws.Range("C13").Value = '1285870 - Marcos Villela Barcza'

$ws.Range("B14").Value = '1.Petróleo: histórico, constituinte, composição e classificação 
2.Geologia do petróleo: origem.
3.Prospecção de petróleo: métodos geológicos, potenciais, sísmicos; 
4.Perfuração: equipamentos, operações;
5.Completação e reservatórios: tipos, etapas, equipamentos, reservatórios;
6.Elevação: elevação natural, sistemas de bombeamentos;
7.Processamento primário: separação do gás natural, tratamento do óleo, tratamento da água, unidade de processamento de gás natural;
8.Derivados do petróleo: tipos, características, gás liquefeito de petróleo, gasolina automotiva, querosene de aviação, óleo diesel, óleos combustíveis industriais, óleos combustíveis marítimos, produtos especiais;
9.Processos de refino: objetivo, tipos de processos, esquemas de refino.
9.1- Destilação de petróleo: equipamentos, esquemas típicos, descrição e variáveis do processo;
9.2- Desasfaltação: carga, descrição e variáveis do processo, produtos;
9.3- Coqueamento retardado: carga, descrição e variáveis do processo, produtos;
9.4- Craqueamento catalítico: carga, descrição e variáveis do processo, produtos;
9.5- Hidrorrefino: carga, descrição e variáveis do processo, produtos;
9.6- Reforma catalítica: carga, descrição e variáveis do processo, produtos;
9.7- Alquilação e isomerização: carga, descrição e variáveis do processo, produtos;
9.8- Tratamento de derivados: tratamento com aminas, tratamentos cáusticos;
9.9- Geração de hidrogênio: carga, descrição e variáveis do processo;
9.10- Recuperação de Enxofre: Processo Claus.
10.Óleos básicos lubrificantes e parafinas: carga, descrição e variáveis do processo.'
$ws.Range("C14").Value = '1.Petróleo: histórico, constituinte, composição e classificação 
2.Geologia do petróleo: origem.
3.Prospecção de petróleo: métodos geológicos, potenciais, sísmicos; 
4.Perfuração: equipamentos, operações;
5.Completação e reservatórios: tipos, etapas, equipamentos, reservatórios;
6.Elevação: elevação natural, sistemas de bombeamentos;
7.Processamento primário: separação do gás natural, tratamento do óleo, tratamento da água, unidade de processamento de gás natural;
8.Derivados do petróleo: tipos, características, gás liquefeito de petróleo, gasolina automotiva, querosene de aviação, óleo diesel, óleos combustíveis industriais, óleos combustíveis marítimos, produtos especiais;
9.Processos de refino: objetivo, tipos de processos, esquemas de refino.
9.1- Destilação de petróleo: equipamentos, esquemas típicos, descrição e variáveis do processo;
9.2- Desasfaltação: carga, descrição e variáveis do processo, produtos;
9.3- Coqueamento retardado: carga, descrição e variáveis do processo, produtos;
9.4- Craqueamento catalítico: carga, descrição e variáveis do processo, produtos;
9.5- Hidrorrefino: carga, descrição e variáveis do processo, produtos;
9.6- Reforma catalítica: carga, descrição e variáveis do processo, produtos;
9.7- Alquilação e isomerização: carga, descrição e variáveis do processo, produtos;
9.8- Tratamento de derivados: tratamento com aminas, tratamentos cáusticos;
9.9- Geração de hidrogênio: carga, descrição e variáveis do processo;
9.10- Recuperação de Enxofre: Processo Claus.
10.Óleos básicos lubrificantes e parafinas: carga, descrição e variáveis do processo.'

$ws.Range("B16").Value = 'Petróleo: histórico, constituinte, composição e classificação; Geologia do petróleo: origem; Prospecção de petróleo: métodos geológicos, potenciais, sísmicos; Perfuração: equipamentos, operações; Completação e reservatórios: tipos, etapas, equipamentos, reservatórios; Elevação: elevação natural, bombeios; Processamento primário: separação do gás natural, tratamento do óleo, tratamento da água, unidade de processamento de gás natural; Derivados do petróleo: tipos, características, gás liquefeito de petróleo, gasolina automotiva, querosene de aviação, óleo diesel, óleos combustíveis industriais, óleos combustíveis marítimos, produtos especiais; Processos de refino: objetivo, tipos de processos, esquemas de refino; Destilação de petróleo: equipamentos, esquemas típicos, descrição e variáveis do processo; Desasfaltação: carga, descrição e variáveis do processo, produtos; Coqueamento retardado: carga, descrição e variáveis do processo, produtos; Craqueamento catalítico: carga, descrição e variáveis do processo, produtos; Hidrorrefino: carga, descrição e variáveis do processo, produtos; Reforma catalítica: carga, descrição e variáveis do processo, produtos; Alquilação e isomerização: carga, descrição e variáveis do processo, produtos; Tratamento de derivados: tratamento com aminas, tratamentos cáusticos; Geração de hidrogênio: carga, descrição e variáveis do processo; Recuperação de Enxofre: Processo Claus; Óleos básicos lubrificantes e parafinas: carga, descrição e variáveis do processo.'
$ws.Range("C16").Value = 'Petróleo: histórico, constituinte, composição e classificação; Geologia do petróleo: origem; Prospecção de petróleo: métodos geológicos, potenciais, sísmicos; Perfuração: equipamentos, operações; Completação e reservatórios: tipos, etapas, equipamentos, reservatórios; Elevação: elevação natural, bombeios; Processamento primário: separação do gás natural, tratamento do óleo, tratamento da água, unidade de processamento de gás natural; Derivados do petróleo: tipos, características, gás liquefeito de petróleo, gasolina automotiva, querosene de aviação, óleo diesel, óleos combustíveis industriais, óleos combustíveis marítimos, produtos especiais; Processos de refino: objetivo, tipos de processos, esquemas de refino; Destilação de petróleo: equipamentos, esquemas típicos, descrição e variáveis do processo; Desasfaltação: carga, descrição e variáveis do processo, produtos; Coqueamento retardado: carga, descrição e variáveis do processo, produtos; Craqueamento catalítico: carga, descrição e variáveis do processo, produtos; Hidrorrefino: carga, descrição e variáveis do processo, produtos; Reforma catalítica: carga, descrição e variáveis do processo, produtos; Alquilação e isomerização: carga, descrição e variáveis do processo, produtos; Tratamento de derivados: tratamento com aminas, tratamentos cáusticos; Geração de hidrogênio: carga, descrição e variáveis do processo; Recuperação de Enxofre: Processo Claus; Óleos básicos lubrificantes e parafinas: carga, descrição e variáveis do processo.'

$ws.Range("B19").Value = 'Aulas expositivas, desenvolvimento de exercícios em sala e fora de sala de aula, discussão de casos práticos e seminários'
$ws.Range("C19").Value = 'Aulas expositivas, desenvolvimento de exercícios em sala e fora de sala de aula, discussão de casos práticos e seminários'

$ws.Range("B20").Value = 'Provas, avaliação através de exercícios ou casos práticos elaborados fora de sala de aula.'
$ws.Range("C20").Value = 'Provas, avaliação através de exercícios ou casos práticos elaborados fora de sala de aula.'

$ws.Range("B21").Value = 'Frequência mínima de 70% e nota igual ou superior a 3,00 e inferior a 5,00 possibilita prova de recuperação.'
$ws.Range("C21").Value = 'Frequência mínima de 70% e nota igual ou superior a 3,00 e inferior a 5,00 possibilita prova de recuperação.'

$ws.Range("B22").Value = 'a)Speight, J. G., The Chemistry and Technology of Petroleum, CRC Press, 4ª Edição, 2007;
b)Thomas, J. E. (Organizador), Fundamentos de Engenharia de Petróleo, Editora Interciência, 2ª Edição, 2004;
c)Brasil, N. I., Araújo, M. A. S., Souza, E. C. M, Processamento de Petróleo e Gás, Editora LTC, 1ª Edição, 2012;
d)Fundamentos do Refino do Petróleo  Tecnologia e Economia, Szklo, A. S., Uller, V. C., Bonfá, M. H. P., Editora Interciência, 3ª Edição, 2012.
e)Oil and Gas Journal;
f)Revista Petro & Química.'
$ws.Range("C22").Value = 'a)Speight, J. G., The Chemistry and Technology of Petroleum, CRC Press, 4ª Edição, 2007;
b)Thomas, J. E. (Organizador), Fundamentos de Engenharia de Petróleo, Editora Interciência, 2ª Edição, 2004;
c)Brasil, N. I., Araújo, M. A. S., Souza, E. C. M, Processamento de Petróleo e Gás, Editora LTC, 1ª Edição, 2012;
d)Fundamentos do Refino do Petróleo  Tecnologia e Economia, Szklo, A. S., Uller, V. C., Bonfá, M. H. P., Editora Interciência, 3ª Edição, 2012.
e)Oil and Gas Journal;
f)Revista Petro & Química.'

# Re-assert the explicit row heights (writing the long wrapped text above can
# trigger an auto-fit pass that nudges a row's height away from the intended,
# previously-set custom height).
$ws.Rows.Item(10).RowHeight = 60
$ws.Rows.Item(11).RowHeight = 60
$ws.Rows.Item(14).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 60
$ws.Rows.Item(16).RowHeight = 120
$ws.Rows.Item(17).RowHeight = 120
$ws.Rows.Item(19).RowHeight = 60
$ws.Rows.Item(20).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 60
$ws.Rows.Item(22).RowHeight = 120
$ws.Rows.Item(24).RowHeight = 30
